# memory-map.xlsx update: "docs: add VideoA Bright Adjust"
#
# 1) Documents a new memory-map entry for addresses 206004/206005
#    ("VideoA Bright/Beight Adjust"), mirroring the existing
#    "VideoA Color Adjust" legend block found at rows 533/534
#    (columns A/D/E/F carry the category name + value/meaning/description
#    legend alongside the main B=address / C=value columns).
# 2) Normalizes the long-standing hex address/byte values in columns B/C
#    (previously mixed-case, e.g. "20010a", "3d") to uppercase
#    (e.g. "20010A", "3D") across the whole sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New legend block for the VideoA Bright/Beight Adjust memory addresses
# ---------------------------------------------------------------------------

# Row 531 -> address 206004
$ws.Range('A531').Value = 'VideoA Beight Adjust'
$ws.Range('D531').NumberFormat = '@'
$ws.Range('D531').Value = '01'
$ws.Range('E531').NumberFormat = '@'
$ws.Range('E531').Value = '01'
$ws.Range('F531').NumberFormat = '@'
$ws.Range('F531').Value = '00 is negative, 01 is positive. See address 206005.'

# Row 532 -> address 206005
$ws.Range('A532').Value = 'VideoA Bright Adjust'
$ws.Range('D532').NumberFormat = '@'
$ws.Range('D532').Value = '44'
$ws.Range('E532').NumberFormat = '@'
$ws.Range('E532').Value = '3D'
$ws.Range('F532').Value = '7-bit signed integer, 0 value is 00 with byte 206004 set to 01. 3Dh=61d 44h=68d (127-59?)'

# ---------------------------------------------------------------------------
# 2) Uppercase the pre-existing hex address / byte-value strings
# ---------------------------------------------------------------------------

$hexUpdates = @(
    @{ Cell = 'C5'; Value = '3D' },
    @{ Cell = 'B12'; Value = '20010A' },
    @{ Cell = 'B13'; Value = '20010B' },
    @{ Cell = 'B14'; Value = '20010C' },
    @{ Cell = 'B15'; Value = '20010D' },
    @{ Cell = 'C15'; Value = '0C' },
    @{ Cell = 'B16'; Value = '20010E' },
    @{ Cell = 'B17'; Value = '20010F' },
    @{ Cell = 'C19'; Value = '5D' },
    @{ Cell = 'C21'; Value = '5F' },
    @{ Cell = 'C23'; Value = '0C' },
    @{ Cell = 'B28'; Value = '20011A' },
    @{ Cell = 'C36'; Value = '3B' },
    @{ Cell = 'B39'; Value = '20020A' },
    @{ Cell = 'B40'; Value = '20020B' },
    @{ Cell = 'C40'; Value = '1B' },
    @{ Cell = 'B41'; Value = '20020C' },
    @{ Cell = 'B42'; Value = '20020D' },
    @{ Cell = 'C42'; Value = '7D' },
    @{ Cell = 'B43'; Value = '20020E' },
    @{ Cell = 'B44'; Value = '20020F' },
    @{ Cell = 'C48'; Value = '5F' },
    @{ Cell = 'B55'; Value = '20021A' },
    @{ Cell = 'C59'; Value = '5C' },
    @{ Cell = 'B66'; Value = '20030A' },
    @{ Cell = 'B67'; Value = '20030B' },
    @{ Cell = 'C67'; Value = '1B' },
    @{ Cell = 'B68'; Value = '20030C' },
    @{ Cell = 'B69'; Value = '20030D' },
    @{ Cell = 'C69'; Value = '7D' },
    @{ Cell = 'B70'; Value = '20030E' },
    @{ Cell = 'B71'; Value = '20030F' },
    @{ Cell = 'C75'; Value = '5F' },
    @{ Cell = 'B82'; Value = '20031A' },
    @{ Cell = 'C86'; Value = '2B' },
    @{ Cell = 'C92'; Value = '4D' },
    @{ Cell = 'B93'; Value = '20040A' },
    @{ Cell = 'B94'; Value = '20040B' },
    @{ Cell = 'C94'; Value = '5B' },
    @{ Cell = 'B95'; Value = '20040C' },
    @{ Cell = 'B96'; Value = '20040D' },
    @{ Cell = 'C96'; Value = '0C' },
    @{ Cell = 'B97'; Value = '20040E' },
    @{ Cell = 'B98'; Value = '20040F' },
    @{ Cell = 'C98'; Value = '0B' },
    @{ Cell = 'C102'; Value = '5F' },
    @{ Cell = 'B109'; Value = '20041A' },
    @{ Cell = 'C115'; Value = '5C' },
    @{ Cell = 'C117'; Value = '0B' },
    @{ Cell = 'B120'; Value = '20050A' },
    @{ Cell = 'B121'; Value = '20050B' },
    @{ Cell = 'B122'; Value = '20050C' },
    @{ Cell = 'B123'; Value = '20050D' },
    @{ Cell = 'B124'; Value = '20050E' },
    @{ Cell = 'B125'; Value = '20050F' },
    @{ Cell = 'C125'; Value = '0D' },
    @{ Cell = 'B136'; Value = '20051A' },
    @{ Cell = 'B147'; Value = '20060A' },
    @{ Cell = 'B148'; Value = '20060B' },
    @{ Cell = 'B149'; Value = '20060C' },
    @{ Cell = 'B150'; Value = '20060D' },
    @{ Cell = 'B151'; Value = '20060E' },
    @{ Cell = 'B152'; Value = '20060F' },
    @{ Cell = 'C154'; Value = '5D' },
    @{ Cell = 'C156'; Value = '5F' },
    @{ Cell = 'C158'; Value = '0F' },
    @{ Cell = 'B163'; Value = '20061A' },
    @{ Cell = 'C169'; Value = '6A' },
    @{ Cell = 'C171'; Value = '0B' },
    @{ Cell = 'B174'; Value = '20070A' },
    @{ Cell = 'B175'; Value = '20070B' },
    @{ Cell = 'B176'; Value = '20070C' },
    @{ Cell = 'B177'; Value = '20070D' },
    @{ Cell = 'C177'; Value = '7D' },
    @{ Cell = 'B178'; Value = '20070E' },
    @{ Cell = 'B179'; Value = '20070F' },
    @{ Cell = 'C179'; Value = '0C' },
    @{ Cell = 'C183'; Value = '5F' },
    @{ Cell = 'C185'; Value = '0D' },
    @{ Cell = 'B190'; Value = '20071A' },
    @{ Cell = 'C192'; Value = '7E' },
    @{ Cell = 'C196'; Value = '6A' },
    @{ Cell = 'B201'; Value = '20100A' },
    @{ Cell = 'B202'; Value = '20100B' },
    @{ Cell = 'B203'; Value = '20100C' },
    @{ Cell = 'B204'; Value = '20100D' },
    @{ Cell = 'C204'; Value = '1F' },
    @{ Cell = 'B205'; Value = '20100E' },
    @{ Cell = 'B206'; Value = '20100F' },
    @{ Cell = 'C218'; Value = '1B' },
    @{ Cell = 'C224'; Value = '2E' },
    @{ Cell = 'B227'; Value = '20110A' },
    @{ Cell = 'B228'; Value = '20110B' },
    @{ Cell = 'B229'; Value = '20110C' },
    @{ Cell = 'B230'; Value = '20110D' },
    @{ Cell = 'C230'; Value = '1F' },
    @{ Cell = 'B231'; Value = '20110E' },
    @{ Cell = 'B232'; Value = '20110F' },
    @{ Cell = 'C246'; Value = '5D' },
    @{ Cell = 'C248'; Value = '5B' },
    @{ Cell = 'C250'; Value = '5A' },
    @{ Cell = 'C252'; Value = '5C' },
    @{ Cell = 'B253'; Value = '20200A' },
    @{ Cell = 'B254'; Value = '20200B' },
    @{ Cell = 'C254'; Value = '5E' },
    @{ Cell = 'B255'; Value = '20200C' },
    @{ Cell = 'B256'; Value = '20200D' },
    @{ Cell = 'C256'; Value = '5F' },
    @{ Cell = 'B257'; Value = '20200E' },
    @{ Cell = 'B258'; Value = '20200F' },
    @{ Cell = 'C265'; Value = '0B' },
    @{ Cell = 'C267'; Value = '0C' },
    @{ Cell = 'B270'; Value = '20210A' },
    @{ Cell = 'B271'; Value = '20210B' },
    @{ Cell = 'C271'; Value = '0F' },
    @{ Cell = 'B272'; Value = '20210C' },
    @{ Cell = 'B273'; Value = '20210D' },
    @{ Cell = 'B274'; Value = '20210E' },
    @{ Cell = 'B275'; Value = '20210F' },
    @{ Cell = 'B287'; Value = '20220A' },
    @{ Cell = 'B288'; Value = '20220B' },
    @{ Cell = 'B289'; Value = '20220C' },
    @{ Cell = 'B290'; Value = '20220D' },
    @{ Cell = 'B291'; Value = '20220E' },
    @{ Cell = 'B292'; Value = '20220F' },
    @{ Cell = 'C297'; Value = '4E' },
    @{ Cell = 'C301'; Value = '2E' },
    @{ Cell = 'C303'; Value = '1A' },
    @{ Cell = 'B304'; Value = '20230A' },
    @{ Cell = 'B305'; Value = '20230B' },
    @{ Cell = 'B306'; Value = '20230C' },
    @{ Cell = 'B307'; Value = '20230D' },
    @{ Cell = 'C307'; Value = '1A' },
    @{ Cell = 'B308'; Value = '20230E' },
    @{ Cell = 'B309'; Value = '20230F' },
    @{ Cell = 'C314'; Value = '3B' },
    @{ Cell = 'B321'; Value = '20300A' },
    @{ Cell = 'B322'; Value = '20300B' },
    @{ Cell = 'B323'; Value = '20300C' },
    @{ Cell = 'B324'; Value = '20300D' },
    @{ Cell = 'B325'; Value = '20300E' },
    @{ Cell = 'B326'; Value = '20300F' },
    @{ Cell = 'C334'; Value = '3B' },
    @{ Cell = 'B341'; Value = '20310A' },
    @{ Cell = 'B342'; Value = '20310B' },
    @{ Cell = 'B343'; Value = '20310C' },
    @{ Cell = 'B344'; Value = '20310D' },
    @{ Cell = 'B345'; Value = '20310E' },
    @{ Cell = 'B346'; Value = '20310F' },
    @{ Cell = 'B361'; Value = '20320A' },
    @{ Cell = 'B362'; Value = '20320B' },
    @{ Cell = 'B363'; Value = '20320C' },
    @{ Cell = 'B364'; Value = '20320D' },
    @{ Cell = 'B365'; Value = '20320E' },
    @{ Cell = 'B366'; Value = '20320F' },
    @{ Cell = 'B381'; Value = '20330A' },
    @{ Cell = 'B382'; Value = '20330B' },
    @{ Cell = 'B383'; Value = '20330C' },
    @{ Cell = 'B384'; Value = '20330D' },
    @{ Cell = 'B385'; Value = '20330E' },
    @{ Cell = 'B386'; Value = '20330F' },
    @{ Cell = 'C394'; Value = '4F' },
    @{ Cell = 'B401'; Value = '20340A' },
    @{ Cell = 'B402'; Value = '20340B' },
    @{ Cell = 'B403'; Value = '20340C' },
    @{ Cell = 'B404'; Value = '20340D' },
    @{ Cell = 'B405'; Value = '20340E' },
    @{ Cell = 'B406'; Value = '20340F' },
    @{ Cell = 'B421'; Value = '20350A' },
    @{ Cell = 'B422'; Value = '20350B' },
    @{ Cell = 'B423'; Value = '20350C' },
    @{ Cell = 'B424'; Value = '20350D' },
    @{ Cell = 'B425'; Value = '20350E' },
    @{ Cell = 'B426'; Value = '20350F' },
    @{ Cell = 'B441'; Value = '20360A' },
    @{ Cell = 'B442'; Value = '20360B' },
    @{ Cell = 'B443'; Value = '20360C' },
    @{ Cell = 'B444'; Value = '20360D' },
    @{ Cell = 'B445'; Value = '20360E' },
    @{ Cell = 'B446'; Value = '20360F' },
    @{ Cell = 'C446'; Value = '0B' },
    @{ Cell = 'B461'; Value = '20370A' },
    @{ Cell = 'B462'; Value = '20370B' },
    @{ Cell = 'B463'; Value = '20370C' },
    @{ Cell = 'B464'; Value = '20370D' },
    @{ Cell = 'B465'; Value = '20370E' },
    @{ Cell = 'B466'; Value = '20370F' },
    @{ Cell = 'C466'; Value = '0B' },
    @{ Cell = 'C468'; Value = '0E' },
    @{ Cell = 'C474'; Value = '0E' },
    @{ Cell = 'C478'; Value = '0E' },
    @{ Cell = 'B481'; Value = '20400A' },
    @{ Cell = 'B482'; Value = '20400B' },
    @{ Cell = 'B483'; Value = '20400C' },
    @{ Cell = 'B484'; Value = '20400D' },
    @{ Cell = 'C484'; Value = '0A' },
    @{ Cell = 'B485'; Value = '20400E' },
    @{ Cell = 'B486'; Value = '20400F' },
    @{ Cell = 'B497'; Value = '20500A' },
    @{ Cell = 'B498'; Value = '20500B' },
    @{ Cell = 'B499'; Value = '20500C' },
    @{ Cell = 'B500'; Value = '20500D' },
    @{ Cell = 'B501'; Value = '20500E' },
    @{ Cell = 'B502'; Value = '20500F' },
    @{ Cell = 'C509'; Value = '0A' },
    @{ Cell = 'B513'; Value = '20501A' },
    @{ Cell = 'B514'; Value = '20501B' },
    @{ Cell = 'B515'; Value = '20501C' },
    @{ Cell = 'B516'; Value = '20501D' },
    @{ Cell = 'B517'; Value = '20501E' },
    @{ Cell = 'B518'; Value = '20501F' },
    @{ Cell = 'E534'; Value = '0B' },
    @{ Cell = 'F534'; Value = '7-bit signed integer, 0 value is 00 with byte 206006 set to 01. Negative values start at 7F.' },
    @{ Cell = 'B537'; Value = '20600A' },
    @{ Cell = 'B538'; Value = '20600B' },
    @{ Cell = 'B539'; Value = '20600C' },
    @{ Cell = 'B540'; Value = '20600D' },
    @{ Cell = 'B541'; Value = '20600E' },
    @{ Cell = 'B542'; Value = '20600F' },
    @{ Cell = 'C542'; Value = '0A' },
    @{ Cell = 'B553'; Value = '20601A' },
    @{ Cell = 'B554'; Value = '20601B' },
    @{ Cell = 'C554'; Value = '0A' },
    @{ Cell = 'B555'; Value = '20601C' },
    @{ Cell = 'B556'; Value = '20601D' },
    @{ Cell = 'B557'; Value = '20601E' },
    @{ Cell = 'B558'; Value = '20601F' },
    @{ Cell = 'B569'; Value = '20602A' },
    @{ Cell = 'B570'; Value = '20602B' },
    @{ Cell = 'B571'; Value = '20602C' },
    @{ Cell = 'B572'; Value = '20602D' },
    @{ Cell = 'B573'; Value = '20602E' },
    @{ Cell = 'B574'; Value = '20602F' },
    @{ Cell = 'B577'; Value = '2F7F00' }
)

foreach ($u in $hexUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---------------------------------------------------------------------------
# View state
# ---------------------------------------------------------------------------
$ws.Range('E537').Select()
